$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-08-28 00:43:04"

$wsZhCn.Range("H4").Value = "2016-08-28 00:42:57"
$wsZhCn.Range("K4").Value = "2016-08-28 00:43:26"

$wsDeDe.Range("H4").Value = "2016-08-28 00:43:04"
$wsDeDe.Range("K4").Value = "2016-08-28 00:43:32"
